$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 57.319636
$ws.Range("H2").Value = 171.958908
$ws.Range("I2").Value = 0.5476981520382651
$ws.Range("J2").Value = 0.5476981520382651
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 8800.227388984356
$ws.Range("R2").Value = 79202.04650085921
$ws.Range("S2").Value = 0.1737411894704793
$ws.Range("T2").Value = 0.1737411894704793
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 57.319636
$ws.Range("H3").Value = 171.958908
$ws.Range("I3").Value = 0.5476981520382651
$ws.Range("J3").Value = 0.5476981520382651
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 9675.537513761565
$ws.Range("R3").Value = 87079.83762385408
$ws.Range("S3").Value = 0.1910222681872302
$ws.Range("T3").Value = 0.1910222681872302
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 57.319636
$ws.Range("H4").Value = 171.958908
$ws.Range("I4").Value = 0.5476981520382651
$ws.Range("J4").Value = 0.5476981520382651
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 3902.912548588974
$ws.Range("R4").Value = 35126.21293730076
$ws.Range("S4").Value = 0.07705444855208086
$ws.Range("T4").Value = 0.07705444855208088
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 57.319636
$ws.Range("H5").Value = 171.958908
$ws.Range("I5").Value = 0.5476981520382651
$ws.Range("J5").Value = 0.5476981520382651
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 5362.978359547028
$ws.Range("R5").Value = 48266.80523592325
$ws.Range("S5").Value = 0.1058802458284748
$ws.Range("T5").Value = 0.1058802458284748
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.975184333333334
$ws.Range("H6").Value = 17.925553
$ws.Range("I6").Value = 0.05709382762749331
$ws.Range("J6").Value = 0.05709382762749331
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 917.3641790821951
$ws.Range("R6").Value = 8256.277611739755
$ws.Range("S6").Value = 0.01811134378764559
$ws.Range("T6").Value = 0.01811134378764559
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.975184333333334
$ws.Range("H7").Value = 17.925553
$ws.Range("I7").Value = 0.05709382762749331
$ws.Range("J7").Value = 0.05709382762749331
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 1008.609338845192
$ws.Range("R7").Value = 9077.484049606726
$ws.Range("S7").Value = 0.01991277935174145
$ws.Range("T7").Value = 0.01991277935174145
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.975184333333334
$ws.Range("H8").Value = 17.925553
$ws.Range("I8").Value = 0.05709382762749331
$ws.Range("J8").Value = 0.05709382762749331
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 406.8522332329345
$ws.Range("R8").Value = 3661.67009909641
$ws.Range("S8").Value = 0.00803240505229365
$ws.Range("T8").Value = 0.008032405052293651
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.975184333333334
$ws.Range("H9").Value = 17.925553
$ws.Range("I9").Value = 0.05709382762749331
$ws.Range("J9").Value = 0.05709382762749331
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 559.0542178943897
$ws.Range("R9").Value = 5031.487961049507
$ws.Range("S9").Value = 0.01103729943581262
$ws.Range("T9").Value = 0.01103729943581262
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 32.32302533333333
$ws.Range("H10").Value = 96.969076
$ws.Range("I10").Value = 0.3088515991858827
$ws.Range("J10").Value = 0.3088515991858827
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 4962.522316667105
$ws.Range("R10").Value = 44662.70085000395
$ws.Range("S10").Value = 0.0979741195268192
$ws.Range("T10").Value = 0.09797411952681921
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 32.32302533333333
$ws.Range("H11").Value = 96.969076
$ws.Range("I11").Value = 0.3088515991858827
$ws.Range("J11").Value = 0.3088515991858827
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 5456.117065553802
$ws.Range("R11").Value = 49105.05358998421
$ws.Range("S11").Value = 0.1077190653103002
$ws.Range("T11").Value = 0.1077190653103002
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 32.32302533333333
$ws.Range("H12").Value = 96.969076
$ws.Range("I12").Value = 0.3088515991858827
$ws.Range("J12").Value = 0.3088515991858827
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 2200.885246058191
$ws.Range("R12").Value = 19807.96721452372
$ws.Range("S12").Value = 0.04345165228535191
$ws.Range("T12").Value = 0.04345165228535192
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 32.32302533333333
$ws.Range("H13").Value = 96.969076
$ws.Range("I13").Value = 0.3088515991858827
$ws.Range("J13").Value = 0.3088515991858827
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 3024.228649633383
$ws.Range("R13").Value = 27218.05784670045
$ws.Range("S13").Value = 0.05970676206341142
$ws.Range("T13").Value = 0.05970676206341142
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.037676333333334
$ws.Range("H14").Value = 27.113029
$ws.Range("I14").Value = 0.08635642114835883
$ws.Range("J14").Value = 0.08635642114835884
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 1387.54556643339
$ws.Range("R14").Value = 12487.91009790051
$ws.Range("S14").Value = 0.02739404409690482
$ws.Range("T14").Value = 0.02739404409690483
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.037676333333334
$ws.Range("H15").Value = 27.113029
$ws.Range("I15").Value = 0.08635642114835883
$ws.Range("J15").Value = 0.08635642114835884
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 1525.557077864237
$ws.Range("R15").Value = 13730.01370077813
$ws.Range("S15").Value = 0.03011877870849325
$ws.Range("T15").Value = 0.03011877870849325
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.037676333333334
$ws.Range("H16").Value = 27.113029
$ws.Range("I16").Value = 0.08635642114835883
$ws.Range("J16").Value = 0.08635642114835884
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 615.3783037186811
$ws.Range("R16").Value = 5538.40473346813
$ws.Range("S16").Value = 0.01214929498256395
$ws.Range("T16").Value = 0.01214929498256396
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.037676333333334
$ws.Range("H17").Value = 27.113029
$ws.Range("I17").Value = 0.08635642114835883
$ws.Range("J17").Value = 0.08635642114835884
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 845.5891554555058
$ws.Range("R17").Value = 7610.302399099552
$ws.Range("S17").Value = 0.01669430336039681
$ws.Range("T17").Value = 0.01669430336039681
